# Remove old parameters for TB NH
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete the 4 rows that hold the now-removed "old tbnh params" group
# (perc_sp_cure_smearpos / lambda_timeto_sp_or_death_smearpos /
#  perc_sp_cure_closed_tb / lambda_timeto_sp_or_death_closed_tb), rows 45-48.
$ws.Range("A45:D48").EntireRow.Delete() | Out-Null

# Recolor the freed-up highlight style (now reused for the trailing two rows)
# from the dark navy tint to a light grey tint.
$ws.Range("A77:D78").Interior.ThemeColor = 3
$ws.Range("A77:D78").Interior.TintAndShade = 0.79998168889431442

# Restore the view state captured in the saved workbook
$ws.Application.ActiveWindow.Zoom = 190
$ws.Range("A31").Select() | Out-Null
$ws.Application.ActiveWindow.ScrollRow = 31
$ws.Range("A49").Select() | Out-Null
